# Auto commit at 2025-11-04  8:31:07.40
# Applies the recorded edit to the workbook:
#  - Metrics!B2:B13 literal values bumped
#  - today!B3:B6 cleared
#  - today!B11:B14 now pull from Metrics via formula (and E11:E14 mirror B11:B14)
#  - selections / active-sheet bookkeeping moved from Chargingdata to today

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metrics sheet: bump the twelve input values
# ---------------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 38765.229999999996
$metrics.Range("B3").Value  = 33901.119999999995
$metrics.Range("B4").Value  = 11945.57
$metrics.Range("B5").Value  = 1608
$metrics.Range("B6").Value  = 4835010.9799999995
$metrics.Range("B7").Value  = 4075977.8000000007
$metrics.Range("B8").Value  = 1418905.4
$metrics.Range("B9").Value  = 187815
$metrics.Range("B10").Value = 33300391.970000003
$metrics.Range("B11").Value = 31351252.959999997
$metrics.Range("B12").Value = 11700627.440000001
$metrics.Range("B13").Value = 1285445

# ---------------------------------------------------------------------------
# 2) today sheet: clear B3:B6 and wire B11:B14 / E11:E14 to Metrics / B
# ---------------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")

$today.Range("B3").ClearContents()
$today.Range("B4").ClearContents()
$today.Range("B5").ClearContents()
$today.Range("B6").ClearContents()

$today.Range("B11").Formula = "=Metrics!B2"
$today.Range("E11").Formula = "=B11"

$today.Range("B12").Formula = "=Metrics!B3"
$today.Range("B13").Formula = "=Metrics!B4"
$today.Range("B14").Formula = "=Metrics!B5"
$today.Range("E12:E14").Formula = "=B12"

# ---------------------------------------------------------------------------
# 3) Selections + active sheet bookkeeping
# ---------------------------------------------------------------------------
$metrics.Range("C8").Select()

$today.Activate()
$today.Range("G9").Select()
